# Update the raw data in columns A and B (rows 1-32) to the new values,
# and shrink both column widths by one character (case 1 data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.26837214234553386
$ws.Cells.Item(1, 2).Value = 0.26812717878335945
$ws.Cells.Item(2, 1).Value = -0.18972571145514827
$ws.Cells.Item(2, 2).Value = 0.18914848632841341
$ws.Cells.Item(3, 1).Value = -0.11259286388350453
$ws.Cells.Item(3, 2).Value = 0.11231591655552364
$ws.Cells.Item(4, 1).Value = -0.10431591661359185
$ws.Cells.Item(4, 2).Value = 0.10379397015684333
$ws.Cells.Item(5, 1).Value = -0.10079397019141467
$ws.Cells.Item(5, 2).Value = 0.09901150203193243
$ws.Cells.Item(6, 1).Value = -0.042959333580771641
$ws.Cells.Item(6, 2).Value = 0.042582239968369962
$ws.Cells.Item(7, 1).Value = -0.032582240051589828
$ws.Cells.Item(7, 2).Value = 0.032498008243710519
$ws.Cells.Item(8, 1).Value = -0.022498008329142394
$ws.Cells.Item(8, 2).Value = 0.022369070020312964
$ws.Cells.Item(9, 1).Value = -0.020369070064696793
$ws.Cells.Item(9, 2).Value = 0.020269830672471567
$ws.Cells.Item(10, 1).Value = -0.018269830719129132
$ws.Cells.Item(10, 2).Value = 0.018264681832329543
$ws.Cells.Item(11, 1).Value = -0.024392161028551129
$ws.Cells.Item(11, 2).Value = 0.024365845796257446
$ws.Cells.Item(12, 1).Value = -0.020865845852119591
$ws.Cells.Item(12, 2).Value = 0.020671058430362432
$ws.Cells.Item(13, 1).Value = -0.017171058489768853
$ws.Cells.Item(13, 2).Value = 0.017082378613092253
$ws.Cells.Item(14, 1).Value = -0.0090823786980429588
$ws.Cells.Item(14, 2).Value = 0.0090534589431268131
$ws.Cells.Item(15, 1).Value = -0.0080534589905774112
$ws.Cells.Item(15, 2).Value = 0.008034786729916199
$ws.Cells.Item(16, 1).Value = -0.0060347867834371627
$ws.Cells.Item(16, 2).Value = 0.0060034776903772347
$ws.Cells.Item(17, 1).Value = -0.0040034777448179071
$ws.Cells.Item(17, 2).Value = 0.0039999999343862669
$ws.Cells.Item(18, 1).Value = -0.016103995643007352
$ws.Cells.Item(18, 2).Value = 0.016091381288802609
$ws.Cells.Item(19, 1).Value = -0.012091381312657301
$ws.Cells.Item(19, 2).Value = 0.012016591536459398
$ws.Cells.Item(20, 1).Value = -0.0080165915620113992
$ws.Cells.Item(20, 2).Value = 0.0080056375077219855
$ws.Cells.Item(21, 1).Value = -0.0040056375335515426
$ws.Cells.Item(21, 2).Value = 0.0039999999739439573
$ws.Cells.Item(22, 1).Value = -0.072555622477494808
$ws.Cells.Item(22, 2).Value = 0.072161731602236046
$ws.Cells.Item(23, 1).Value = -0.040495494420451195
$ws.Cells.Item(23, 2).Value = 0.040098312120139035
$ws.Cells.Item(24, 1).Value = -0.020098312245853123
$ws.Cells.Item(24, 2).Value = 0.019999999872634788
$ws.Cells.Item(25, 1).Value = -0.054052168489812757
$ws.Cells.Item(25, 2).Value = 0.053994236069327073
$ws.Cells.Item(26, 1).Value = -0.035358442507927634
$ws.Cells.Item(26, 2).Value = 0.035341363892202438
$ws.Cells.Item(27, 1).Value = -0.032841363935138368
$ws.Cells.Item(27, 2).Value = 0.032750490758456152
$ws.Cells.Item(28, 1).Value = -0.030750490802283537
$ws.Cells.Item(28, 2).Value = 0.030697866021226261
$ws.Cells.Item(29, 1).Value = -0.023697866094749109
$ws.Cells.Item(29, 2).Value = 0.023688785127723655
$ws.Cells.Item(30, 1).Value = 0.036311214509289425
$ws.Cells.Item(30, 2).Value = -0.036416259248358784
$ws.Cells.Item(31, 1).Value = -0.014022148155421021
$ws.Cells.Item(31, 2).Value = 0.014001192772850501
$ws.Cells.Item(32, 1).Value = -0.0040011928608052472
$ws.Cells.Item(32, 2).Value = 0.0039999999446038714

# Column widths shrink by ~1 character (16.43 -> 15.43, 15.71 -> 14.71 in
# stored XML units); ColumnWidth is specified in character-width units.
$ws.Columns(1).ColumnWidth = 14.666666666666666
$ws.Columns(2).ColumnWidth = 13.833333333333332
